# Backend answer-key correction + new frontend question rows
# Applies to the "Questions" worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions")

# ---------------------------------------------------------------
# Row 6 - replace the old joke "creator" question with a real one
# ---------------------------------------------------------------
$ws.Range("A6").Value = "It is called the river of life"
$ws.Range("B6").Value = "Nile River"
$ws.Range("C6").Value = "Blood"
$ws.Range("D6").Value = "Pasig River"
$ws.Range("E6").Value = "Rivermaya"
$ws.Range("F6").Value = "Egyptian River"

# ---------------------------------------------------------------
# Row 7 - replace the old joke "date today" question with a real one
# ---------------------------------------------------------------
$ws.Range("A7").Value = "When is the World War I"
$ws.Range("B7").Value = "Jan. 11, 2023"
$ws.Range("C7").Value = 5323
$ws.Range("C7").NumberFormat = "d-mmm-yy"
$ws.Range("D7").NumberFormat = "d-mmm-yy"
$ws.Range("D7").Value = "Dec. 11, 1998"
$ws.Range("E7").Value = 38785
$ws.Range("F7").Value = "Sep. 06, 1921"

# ---------------------------------------------------------------
# Row 8 (new) - Chemistry question
# ---------------------------------------------------------------
$ws.Range("A8").Value = "What is the chemical symbol for element Gold?"
$ws.Range("B8").Value = "Au"
$ws.Range("C8").Value = "Ag"
$ws.Range("D8").Value = "Fe"
$ws.Range("E8").Value = "Hg"
$ws.Range("F8").Value = "O"
$ws.Range("B8").Interior.Color = 65535

# ---------------------------------------------------------------
# Row 9 (new) - Geography question
# ---------------------------------------------------------------
$ws.Range("A9").Value = "What is the capital of France?"
$ws.Range("B9").Value = "Madrid"
$ws.Range("C9").Value = "Berlin"
$ws.Range("D9").Value = "Paris"
$ws.Range("E9").Value = "Rome"
$ws.Range("F9").Value = "Manila"
$ws.Range("D9").Interior.Color = 65535

# ---------------------------------------------------------------
# Row 10 (new) - Biology question
# ---------------------------------------------------------------
$ws.Range("A10").Value = "What is the powerhouse of the cell?"
$ws.Range("B10").Value = "Nucleus"
$ws.Range("C10").Value = "Mitochondria"
$ws.Range("D10").Value = "Endoplasmic Reticulum"
$ws.Range("E10").Value = "Golgi Apparatus"
$ws.Range("F10").Value = "None of the Above"
$ws.Range("C10").Interior.Color = 65535

# ---------------------------------------------------------------
# Update the view / selection to match the final saved state
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("A8").Select()
